# BAU Vehicle Subsidy.xlsx - "Run csv export tool" commit
# 1) Zero out the Commercial Vehicles credit input cells (B4, B5) - this
#    cascades through every formula that depends on them.
# 2) Refresh the "battery electric vehicles" / "plugin hybrid vehicles"
#    rows (B3:AE3, B4:AE4) on "Passenger Vehicle Calculations" with newly
#    exported data - this also cascades downstream.
# 3) Re-point the active sheet / selection to match where the user ended up
#    (Passenger Vehicle Calculations tab active, selection on O26) and move
#    the "Commercial Vehicles" sheet's remembered selection to B6.

$wb = $excel.ActiveWorkbook

# --- 1) Commercial Vehicles credit inputs -> 0 ---------------------------
$cv = $wb.Worksheets.Item("Commercial Vehicles")
$cv.Range("B4").Value = 0
$cv.Range("B5").Value = 0

# --- 2) Passenger Vehicle Calculations refreshed export rows -------------
$pvc = $wb.Worksheets.Item("Passenger Vehicle Calculations")

$row3 = @(1646.4981124797289, 1368.5857570401329, 983.6971469156598, 654.32041202251673, 219.25443921462275, 29.671981569503362, 29.671981569503377, 29.671981569503377, 29.67198156950337, 29.671981569503377, 29.67198156950338, 29.662930818658605, 29.657500368151737, 29.657500368151744, 20.04560297100296, 14.278464532713683, 14.278464532713679, 14.278464532713679, 14.278464532713675, 14.278464532713675, 14.278464532713679, 14.278464532713674, 14.278464532713679, 14.278464532713683, 14.278464532713679, 14.278464532713683, 14.278464532713679, 14.278464532713679, 10.395692420306119, 8.066029152861578)
$row4 = @(6652.7185168992282, 3617.050692978903, 1422.8620974282658, 838.46462929681411, 421.08837128836086, 300.52351835174784, 300.52351835174778, 300.52351835174778, 300.52351835174778, 300.52351835174778, 300.52351835174778, 245.24790619147615, 149.54513499056355, 58.275878749437524, 11.623644285733851, 2.9811566108910497, 2.9811566108910501, 2.981156610891051, 2.9811566108910506, 2.9811566108910497, 2.9811566108910506, 2.9811566108910497, 2.9811566108910501, 2.9811566108910501, 2.9811566108910497, 2.9811566108910501, 2.9811566108910501, 2.9811566108910501, 1.1179337290841438, 0)

$col = 2
foreach ($v in $row3) {
    $pvc.Cells.Item(3, $col).Value = $v
    $col = $col + 1
}

$col = 2
foreach ($v in $row4) {
    $pvc.Cells.Item(4, $col).Value = $v
    $col = $col + 1
}

# --- 3) View state: active tab + selections -------------------------------
$about = $wb.Worksheets.Item("About")

$pvc.Activate()
$pvc.Range("O26").Select()

$cv.Range("B6").Select()

$about.Range("B36").Select()
$pvc.Activate()
